$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.351.26"
$ws.Range("E2").Value = "  +2.28%  "
$ws.Range("D3").Value = "1.660.85"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").Value = "'220.16"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("E8").Value = "  +1.26%  "
$ws.Range("D9").Value = "'0.0626"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").Value = "'20.07"
$ws.Range("E10").Value = "  +4.88%  "
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").Value = "1.894.16"
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("D13").Value = "1.663.81"
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").Value = "'4.19"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").Value = "'0.533"
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").Value = "'67.06"
$ws.Range("E16").Value = "  +3.64%  "
$ws.Range("D17").Value = "27.338.26"
$ws.Range("E17").Value = "  +2.25%  "
$ws.Range("D18").Value = "0.0₃0736"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "'222.70"
$ws.Range("E19").Value = "  +4.02%  "
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("E21").Value = "  +8.55%  "
$ws.Range("E22").Value = "  +1.56%  "
$ws.Range("E23").Value = "  +5.86%  "
$ws.Range("D24").Value = "'9.29"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'147.03"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("D27").Value = "'7.42"
$ws.Range("E27").Value = "  +3.68%  "
$ws.Range("D28").Value = "'0.119"
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("D29").Value = "'16.05"
$ws.Range("E29").Value = "  +2.38%  "
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("D32").Value = "'3.39"
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("D33").Value = "'3.01"
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("E34").Value = "  +2.28%  "
$ws.Range("D35").Value = "1.264.71"
$ws.Range("E35").Value = "  -1.76%  "
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("D37").Value = "'0.0178"
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("D38").Value = "'0.538"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").Value = "'0.838"
$ws.Range("E39").Value = "  +2.49%  "
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("E41").Value = "  +1.65%  "
$ws.Range("D42").Value = "'5.38"
$ws.Range("E42").Value = "  +1.83%  "
$ws.Range("D43").Value = "1.807.22"
$ws.Range("E43").Value = "  +1.54%  "
$ws.Range("D44").Value = "'2.12"
$ws.Range("E44").Value = "  -4.33%  "
$ws.Range("D45").Value = "'61.82"
$ws.Range("E45").Value = "  +1.24%  "
$ws.Range("D46").Value = "'92.07"
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("E47").Value = "  +1.39%  "
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").Value = "'0.0984"
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("D50").Value = "'7.69"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("E51").Value = "  +0.28%  "
